$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we are about to rewrite to be stored as
# text (matching the source feed, which uses "." as a thousands separator),
# since several new values would otherwise auto-parse as numbers.
$ws.Range("D2:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.290.05'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '1.831.74'
$ws.Range("E3").Value = '  -0.57%  '

$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.59%  '

$ws.Range("D5").Value = '235.01'
$ws.Range("E5").Value = '  -1.90%  '

$ws.Range("D6").Value = '0.6038'
$ws.Range("E6").Value = '  -3.73%  '

$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.27%  '

$ws.Range("D8").Value = '0.07077'
$ws.Range("E8").Value = '  -5.05%  '

$ws.Range("D9").Value = '0.2802'
$ws.Range("E9").Value = '  -3.43%  '

$ws.Range("D10").Value = '23.57'
$ws.Range("E10").Value = '  -5.08%  '

$ws.Range("D11").Value = '0.07662'
$ws.Range("E11").Value = '  -0.73%  '

$ws.Range("D12").Value = '1.824.99'
$ws.Range("E12").Value = '  -0.76%  '

$ws.Range("D13").Value = '4.800'
$ws.Range("E13").Value = '  -3.48%  '

$ws.Range("D14").Value = '0.000009966'
$ws.Range("E14").Value = '  -3.05%  '

$ws.Range("D15").Value = '0.6284'
$ws.Range("E15").Value = '  -7.20%  '

$ws.Range("D16").Value = '2.071.32'
$ws.Range("E16").Value = '  -1.01%  '

$ws.Range("D17").Value = '79.23'
$ws.Range("E17").Value = '  -3.21%  '

$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '5.864'
$ws.Range("E18").Value = '  -6.23%  '

$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '29.205.96'
$ws.Range("E19").Value = '  -0.65%  '

$ws.Range("D20").Value = '226.76'
$ws.Range("E20").Value = '  -2.69%  '

$ws.Range("D21").Value = '1.004'
$ws.Range("E21").Value = '  +0.34%  '

$ws.Range("D22").Value = '11.72'
$ws.Range("E22").Value = '  -4.81%  '

$ws.Range("D23").Value = '7.014'
$ws.Range("E23").Value = '  -4.38%  '

$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  +0.28%  '

$ws.Range("D25").Value = '155.83'
$ws.Range("E25").Value = '  -1.48%  '

$ws.Range("D26").Value = '8.036'
$ws.Range("E26").Value = '  -5.42%  '

$ws.Range("D27").Value = '0.1302'
$ws.Range("E27").Value = '  -3.71%  '

$ws.Range("D28").Value = '16.60'
$ws.Range("E28").Value = '  -4.44%  '

$ws.Range("D29").Value = '1.485'
$ws.Range("E29").Value = '  +1.74%  '

$ws.Range("D30").Value = '0.06206'
$ws.Range("E30").Value = '  -14.60%  '

$ws.Range("D31").Value = '1.444'
$ws.Range("E31").Value = '  -2.49%  '

$ws.Range("D32").Value = '3.836'
$ws.Range("E32").Value = '  -5.22%  '

$ws.Range("D33").Value = '3.805'
$ws.Range("E33").Value = '  -6.28%  '

$ws.Range("D34").Value = '1.125'
$ws.Range("E34").Value = '  -1.36%  '

$ws.Range("D35").Value = '1.751'
$ws.Range("E35").Value = '  -3.77%  '

$ws.Range("D36").Value = '0.6423'
$ws.Range("E36").Value = '  -7.99%  '

$ws.Range("D37").Value = '2.545'
$ws.Range("E37").Value = '  -1.08%  '

$ws.Range("D38").Value = '1.222.40'
$ws.Range("E38").Value = '  -1.02%  '

$ws.Range("D39").Value = '2.739'
$ws.Range("E39").Value = '  -2.73%  '

$ws.Range("D40").Value = '0.01750'
$ws.Range("E40").Value = '  -4.80%  '

$ws.Range("D41").Value = '6.552'
$ws.Range("E41").Value = '  -6.31%  '

$ws.Range("D42").Value = '0.9064'
$ws.Range("E42").Value = '  -4.18%  '

$ws.Range("D43").Value = '1.004'
$ws.Range("E43").Value = '  +0.39%  '

$ws.Range("D44").Value = '1.985.60'
$ws.Range("E44").Value = '  -0.78%  '

$ws.Range("D45").Value = '100.69'
$ws.Range("E45").Value = '  -0.21%  '

$ws.Range("D46").Value = '62.86'
$ws.Range("E46").Value = '  -4.07%  '

$ws.Range("E47").Value = '  -2.00%  '

$ws.Range("D48").Value = '8.525'
$ws.Range("E48").Value = '  -4.44%  '

$ws.Range("D49").Value = '1.585'
$ws.Range("E49").Value = '  -8.09%  '

$ws.Range("D50").Value = '0.4571'
$ws.Range("E50").Value = '  -0.29%  '

$ws.Range("D51").Value = '0.05518'
$ws.Range("E51").Value = '  -2.45%  '
